$wb = $excel.ActiveWorkbook

# --- 1) Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# --- 3) Header row ---
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Match header style (bold, centered, bordered) used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4) Data rows ---
$data = @(
    @(44934.99999999999,46,-319.9685429563868,376.3912971710783),
    @(44997.99999999999,84,-280.4913563981954,422.0849427592263),
    @(45011.99999999999,93,-254.5472908263584,467.760426414392),
    @(45018.99999999999,97,-262.7022757336387,461.9700942561909),
    @(45025.99999999999,101,-248.8019979955205,455.7163132605318),
    @(45039.99999999999,110,-243.7722311638711,480.4463933518363),
    @(45046.99999999999,114,-242.1392877080404,487.7316379006663),
    @(45053.99999999999,119,-245.1690488310549,464.3645426627177),
    @(45067.99999999999,127,-223.5433889681864,485.5430194391527),
    @(45074.99999999999,131,-231.4177988667169,463.8325462909695),
    @(45081.99999999999,136,-220.9405940020765,485.9590137286376),
    @(45088.99999999999,140,-194.3802175659778,519.7215039043109),
    @(45095.99999999999,144,-236.4503977200042,499.5248360768907),
    @(45102.99999999999,149,-208.9341704978302,482.488876162671),
    @(45123.99999999999,161,-163.7650152830959,504.8496375912809),
    @(45137.99999999999,170,-176.4187437450903,521.4327558069159),
    @(45144.99999999999,174,-165.4560269531928,520.8633512875024),
    @(45165.99999999999,187,-164.3556267893002,511.1482257506889),
    @(45172.99999999999,191,-160.3312222758042,560.3707468799902),
    @(45186.99999999999,200,-163.5997568239636,571.7037528985917),
    @(45193.99999999999,204,-152.562261784916,543.2542539908403),
    @(45200.99999999999,208,-138.8786468559357,563.973026838199),
    @(45207.99999999999,213,-162.0186614208338,555.3528346198947),
    @(45214.99999999999,217,-140.671407294941,567.6870920461408),
    @(45221.99999999999,221,-126.9517450539723,574.6347581372793),
    @(45228.99999999999,226,-129.1792556616787,563.41822418098),
    @(45235.99999999999,230,-137.0682621265865,598.1463291197176),
    @(45249.99999999999,238,-110.6685463449515,624.561850143986),
    @(45256.99999999999,243,-108.5245446585746,613.4095217955287),
    @(45263.99999999999,247,-117.2013943398303,614.5563853063397),
    @(45298.99999999999,268,-47.07288024579269,620.0035381107182),
    @(45312.99999999999,277,-100.0203697223883,650.4613110151421),
    @(45333.99999999999,290,-56.84760311385759,637.3905263620478),
    @(45340.99999999999,294,-68.3283101225874,631.7679691118626),
    @(45347.99999999999,298,-63.31436358479044,672.2558479197266),
    @(45354.99999999999,303,-31.54824706818674,668.4462764618791),
    @(45361.99999999999,307,-78.14015460875696,649.1865393647207),
    @(45368.99999999999,311,-36.6470723068541,658.4779213020892),
    @(45375.99999999999,316,18.52022725296345,680.1132854546338),
    @(45382.99999999999,320,-37.65822578241712,674.6064769589825),
    @(45389.99999999999,324,-26.90983956096409,673.9870990725274),
    @(45396.99999999999,328,-19.02759872571339,673.1548667518067),
    @(45410.99999999999,337,-24.04571053562657,696.2396033680282),
    @(45417.99999999999,341,-9.502454328585516,686.0661627029971),
    @(45424.99999999999,346,2.790800104704386,699.4494248521837),
    @(45431.99999999999,350,3.406951863223717,721.049617922164),
    @(45438.99999999999,354,21.37068329864386,720.7517613469406),
    @(45445.99999999999,358,-3.00584603076056,681.4798823409446),
    @(45452.99999999999,363,22.98789688180685,718.7119626088013),
    @(45459.99999999999,367,49.64962653783871,753.1257679660774),
    @(45466.99999999999,371,18.8366155685107,735.427495086809),
    @(45480.99999999999,380,7.745948451030511,724.2844332182709),
    @(45487.99999999999,384,31.67666876272718,739.0014283711428),
    @(45508.99999999999,397,30.19646188130806,770.1782128157218),
    @(45515.99999999999,401,95.67561726391304,763.0361617298474),
    @(45529.99999999999,410,66.40129309108843,765.6063429911634),
    @(45536.99999999999,414,56.86532764708802,779.1958964243117),
    @(45543.99999999999,418,74.3969612762633,801.0576555406134),
    @(45557.99999999999,427,63.57720631138834,819.1841657180864),
    @(45564.99999999999,431,76.43795407557249,774.8071006256999),
    @(45599.99999999999,453,94.37100758435501,820.116707537319),
    @(45606.99999999999,457,104.5635758890956,801.4720767894507),
    @(45613.99999999999,461,121.1114224617793,815.4980317036354),
    @(45620.99999999999,465,124.4861507962638,814.182321186645),
    @(45627.99999999999,470,135.8296722432809,825.036094236673),
    @(45634.99999999999,474,136.8442081189297,830.8948156271952),
    @(45641.99999999999,478,132.5873921572982,841.4662607242471),
    @(45648.99999999999,483,115.6063781267159,839.7292026132671),
    @(45655.99999999999,487,137.74517546868,842.4179404422808)
)

$arr = New-Object 'object[,]' $data.Count,4
for ($i = 0; $i -lt $data.Count; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}
$newSheet.Range("A2:D70").Value2 = $arr

# Match the date-formatted style used for column A on the other sheets
# (copy a single formatted source cell and let PasteSpecial tile it across the target range)
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A" + (1 + $data.Count)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore original active sheet/selection (the new sheet was selected while building it)
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
